# svgMultiple-template.docx: convert the two M2Doc query fields (Word
# fldChar/instrText fields) into plain, brace-delimited literal text runs
# - i.e. replace the `{ FIELD }` constructs with literal "{...}" text so the
# new TokenIteratorFieldRewriterSplit parser can tokenize them straight out
# of the paragraph text.
#
# For each field this:
#   1. reproduces the exact original w:instrText run split as plain w:t runs
#   2. turns the leading/trailing " " delimiter of the field code into the
#      literal "{" / "}" characters (M2Doc's template-tag delimiters)
#   3. removes the w:fldChar begin/end runs entirely
#
# Field 1 (paragraph 2): " " | "m" | ":" | "'" | "doc.html" | "'.fromHTMLURI()" | " "
#      -> "{" | "m" | ":" | "'" | "doc.html" | "'.fromHTMLURI()" | "}"
# Field 2 (paragraph 4): " m:'doc" | "1" | ".html'.fromHTMLURI() "
#      -> "{m:'doc" | "1" | ".html'.fromHTMLURI()}"

$d = $word.ActiveDocument

# Locate the Word.Paragraph that contains document position $pos.
function Get-ParaForPosition($pos) {
    foreach ($p in $d.Paragraphs) {
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $p
        }
    }
    return $null
}

# Delete $field (its begin/instrText/end runs disappear entirely) and insert
# literal text runs - one run per element of $texts - at the start of the
# field's own paragraph, preserving that paragraph's other content/props.
function Replace-FieldWithRuns($field, $texts) {
    $para = Get-ParaForPosition($field.Code.Start)
    $paraStart = $para.Range.Start

    $runXml = ""
    for ($i = 0; $i -lt $texts.Count; $i++) {
        $escaped = $texts[$i] -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
        if ($i -eq $texts.Count - 1) {
            # Match the source formatting: only the final run of the group
            # carries an explicit xml:space="preserve".
            $runXml += "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
        }
        else {
            $runXml += "<w:r><w:t>$escaped</w:t></w:r>"
        }
    }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $field.Delete()
    $insertRange = $d.Range($paraStart, $paraStart)
    $insertRange.InsertXML($xml)
}

# Fields.Item(1) always refers to the first remaining field, so processing
# them in document order and always grabbing Item(1) walks through both.
Replace-FieldWithRuns $d.Fields.Item(1) @("{", "m", ":", "'", "doc.html", "'.fromHTMLURI()", "}")
Replace-FieldWithRuns $d.Fields.Item(1) @("{m:'doc", "1", ".html'.fromHTMLURI()}")
